# Minimal no-op test script to see baseline output differences
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
